# Updates the cryptos list (prices, volume %, and a couple of reordered
# coin rows) on worksheet 1, per the Oct 25 2023 GitHub Actions data refresh.
# Column D holds price text that looks numeric (e.g. "0.997", "1.00"); force
# those specific cells to Text format first so Excel doesn't silently coerce
# them into real numbers (which would drop meaningful trailing zeros, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.804.96"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "1.803.39"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.69"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.558"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.15"
$ws.Range("E8").Value = "  +5.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.286"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0677"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "2.055.99"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.33"
$ws.Range("E13").Value = "  +13.17%  "
$ws.Range("D14").Value = "1.816.44"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "34.761.17"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.98"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "258.25"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").Value = "0.0₃0776"
$ws.Range("E20").Value = "  +4.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.54"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.27"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.72"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.64"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.18"
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.83"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0523"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.63"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.92"
$ws.Range("E34").Value = "  +9.59%  "
$ws.Range("D35").Value = "1.461.77"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0192"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.639"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.84"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.911"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.11"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.06"
$ws.Range("E44").Value = "  +5.37%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0508"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.955.43"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.05"
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.16"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.998"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.70"
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.90"
$ws.Range("E51").Value = "  -1.99%  "
